# Review feedback: the sheet used a literal "Key" / "Value" header row above
# the Table2 data. Per review, that label row is removed entirely (the table
# now starts right at the first real record - "Universal audience") and the
# table is switched to have no header row, since there's no longer a
# dedicated caption row for it to bind to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table's header row is just literal text ("Key"/"Value") sitting in
# worksheet row 1. Turn off the table's header row *before* deleting the
# worksheet row so the table keeps its column definitions (and the same
# number of visible rows: A1:B13 -> A1:B12) instead of trying to re-derive
# headers from the new first row.
$lo = $ws.ListObjects.Item(1)
$lo.ShowHeaders = $false

# Now remove the old "Key"/"Value" row from the sheet - everything below
# shifts up by one (dimension A1:E32 -> A1:E31, the "Manufacturing" row
# drops off the bottom of the table, and the spacer cells in column E
# shift from the even rows to the odd rows).
$ws.Rows(1).Delete()

# Reflect the reviewer's cursor position in the saved file.
$ws.Range("C4").Select() | Out-Null
